$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated 95th-percentile / mean / related stats; river data update May 2024) ---
$ws.Range("G9").Value = 0.207851733930163
$ws.Range("G10").Value = 0.207851733930163
$ws.Range("G19").Value = 0.19329365963562
$ws.Range("G20").Value = 0.19329365963562
$ws.Range("G29").Value = 0.196515028224602
$ws.Range("G30").Value = 0.196515028224602
$ws.Range("G43").Value = 0.178723054524102
$ws.Range("G44").Value = 0.178723054524102
$ws.Range("G57").Value = 0.168697192455137
$ws.Range("G58").Value = 0.168697192455137
$ws.Range("G71").Value = 0.146743744179275
$ws.Range("I71").Value = 0.3395
$ws.Range("L71").Value = 0.07575
$ws.Range("G72").Value = 0.146743744179275
$ws.Range("I72").Value = 0.3395
$ws.Range("L72").Value = 0.07575
$ws.Range("F85").Value = 0.1125
$ws.Range("G85").Value = 0.140420410179753
$ws.Range("N85").Value = 0.34535
$ws.Range("F86").Value = 0.1125
$ws.Range("G86").Value = 0.140420410179753
$ws.Range("N86").Value = 0.34535
$ws.Range("G99").Value = 0.128154905390605
$ws.Range("N99").Value = 0.34535
$ws.Range("G100").Value = 0.128154905390605
$ws.Range("N100").Value = 0.34535
$ws.Range("F113").Value = 0.10875
$ws.Range("G113").Value = 0.136676482014987
$ws.Range("N113").Value = 0.35055
$ws.Range("F114").Value = 0.10875
$ws.Range("G114").Value = 0.136676482014987
$ws.Range("N114").Value = 0.35055
$ws.Range("G131").Value = 0.133226482014987
$ws.Range("N131").Value = 0.35055
$ws.Range("G132").Value = 0.133226482014987
$ws.Range("N132").Value = 0.35055
$ws.Range("G149").Value = 0.11915981534832
$ws.Range("G150").Value = 0.11915981534832
$ws.Range("G156").Value = 1.13132287355318
$ws.Range("H156").Value = 3.34437855025247
$ws.Range("I156").Value = 2.905
$ws.Range("G169").Value = 0.104948333333333
$ws.Range("L169").Value = 0.01495
$ws.Range("G170").Value = 0.104948333333333
$ws.Range("L170").Value = 0.01495
$ws.Range("G176").Value = 1.28143855061721
$ws.Range("I176").Value = 3
$ws.Range("N176").Value = 2.6276
$ws.Range("G189").Value = 0.097523939206117
$ws.Range("I189").Value = 0.3216
$ws.Range("L189").Value = 0.01495
$ws.Range("G190").Value = 0.097523939206117
$ws.Range("I190").Value = 0.3216
$ws.Range("L190").Value = 0.01495
$ws.Range("G196").Value = 1.35217889105919
$ws.Range("I196").Value = 2.905
$ws.Range("F209").Value = 0.04345
$ws.Range("G209").Value = 0.09104045280778481
$ws.Range("I209").Value = 0.3361
$ws.Range("M209").Value = 0.16934
$ws.Range("F210").Value = 0.04345
$ws.Range("G210").Value = 0.09104045280778481
$ws.Range("I210").Value = 0.3361
$ws.Range("M210").Value = 0.16934
$ws.Range("G216").Value = 1.41824565625444
$ws.Range("H216").Value = 4.33713254374192
$ws.Range("I216").Value = 3.25828
$ws.Range("N216").Value = 2.658
$ws.Range("G229").Value = 0.0872049281364904
$ws.Range("I229").Value = 0.33999
$ws.Range("L229").Value = 0.01187
$ws.Range("M229").Value = 0.1665
$ws.Range("G230").Value = 0.0872049281364904
$ws.Range("I230").Value = 0.33999
$ws.Range("L230").Value = 0.01187
$ws.Range("M230").Value = 0.1665
$ws.Range("G236").Value = 1.31629266215084
$ws.Range("H236").Value = 4.33713254374192
$ws.Range("I236").Value = 3.29272
$ws.Range("N236").Value = 2.7188
$ws.Range("F249").Value = 0.06594999999999999
$ws.Range("G249").Value = 0.101575702759533
$ws.Range("I249").Value = 0.3843
$ws.Range("L249").Value = 0.01403
$ws.Range("N249").Value = 0.29129
$ws.Range("F250").Value = 0.06594999999999999
$ws.Range("G250").Value = 0.101575702759533
$ws.Range("I250").Value = 0.3843
$ws.Range("L250").Value = 0.01403
$ws.Range("N250").Value = 0.29129
$ws.Range("G256").Value = 1.16940836196924
$ws.Range("H256").Value = 4.33713254374192
$ws.Range("G269").Value = 0.127444021830895
$ws.Range("I269").Value = 0.4262
$ws.Range("G270").Value = 0.127444021830895
$ws.Range("I270").Value = 0.4262
$ws.Range("G276").Value = 1.08857502863591
$ws.Range("H276").Value = 4.33713254374192
$ws.Range("G289").Value = 0.138380998488156
$ws.Range("G290").Value = 0.138380998488156
$ws.Range("G296").Value = 1.10669184629708
$ws.Range("H296").Value = 4.33713254374192
$ws.Range("G309").Value = 0.130468878907453
$ws.Range("G310").Value = 0.130468878907453
$ws.Range("G329").Value = 0.127795343257627
$ws.Range("G330").Value = 0.127795343257627

# --- Append new rows 335-337 (ASPM, MCI, QMCI for period 2019 - 2023) ---
$ws.Range("A335").Value = "Hautapu at US Rangitikei River Conf"
$ws.Range("B335").Value = "ASPM"
$ws.Range("C335").Value = "D"
$ws.Range("D335").Value = "2019 - 2023"
$ws.Range("E335").Value = "RepSite"
$ws.Range("F335").Value = 0.222
$ws.Range("G335").Value = 0.2594
$ws.Range("H335").Value = 0.42
$ws.Range("I335").Value = 0.42
$ws.Range("L335").Value = 0.251
$ws.Range("M335").Value = 0.3857
$ws.Range("N335").Value = 0.42
$ws.Range("O335").Value = 1842978.43
$ws.Range("P335").Value = 5595723.71
$ws.Range("Q335").Value = "Rangitikei District"
$ws.Range("R335").Value = "Rangitīkei-Turakina"
$ws.Range("S335").Value = "Middle Rangitikei"
$ws.Range("T335").Value = "Rang_2g"

$ws.Range("A336").Value = "Hautapu at US Rangitikei River Conf"
$ws.Range("B336").Value = "MCI"
$ws.Range("C336").Value = "D"
$ws.Range("D336").Value = "2019 - 2023"
$ws.Range("E336").Value = "RepSite"
$ws.Range("F336").Value = 84.44
$ws.Range("G336").Value = 86.3
$ws.Range("H336").Value = 113
$ws.Range("I336").Value = 113
$ws.Range("L336").Value = 82.655
$ws.Range("M336").Value = 106
$ws.Range("N336").Value = 113
$ws.Range("O336").Value = 1842978.43
$ws.Range("P336").Value = 5595723.71
$ws.Range("Q336").Value = "Rangitikei District"
$ws.Range("R336").Value = "Rangitīkei-Turakina"
$ws.Range("S336").Value = "Middle Rangitikei"
$ws.Range("T336").Value = "Rang_2g"

$ws.Range("A337").Value = "Hautapu at US Rangitikei River Conf"
$ws.Range("B337").Value = "QMCI"
$ws.Range("C337").Value = "D"
$ws.Range("D337").Value = "2019 - 2023"
$ws.Range("E337").Value = "RepSite"
$ws.Range("F337").Value = 3.637
$ws.Range("G337").Value = 4.1404
$ws.Range("H337").Value = 5.71
$ws.Range("I337").Value = 5.71
$ws.Range("L337").Value = 3.9375
$ws.Range("M337").Value = 5.20075
$ws.Range("N337").Value = 5.71
$ws.Range("O337").Value = 1842978.43
$ws.Range("P337").Value = 5595723.71
$ws.Range("Q337").Value = "Rangitikei District"
$ws.Range("R337").Value = "Rangitīkei-Turakina"
$ws.Range("S337").Value = "Middle Rangitikei"
$ws.Range("T337").Value = "Rang_2g"

